# "instructions for FC and Prod"
#
# The original workbook's column-I formulas compared the trigger's Number
# (column H) against the literal "singular" - but the actual data in
# column H only ever contains the short codes "sg" / "pl", so every row
# evaluated to the "plural" instruction regardless of the real value.
# This fixes the comparison to use "sg" (matching the data), which in turn
# flips the cached result to "singular" wherever H is actually "sg".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite I2:I31 with the corrected IF() test. Re-entering the same formula
# text on every row lets the shared-formula group (si=0) be rebuilt
# naturally on save, same as it was before the edit.
for ($r = 2; $r -le 31; $r++) {
    $formula = '=IF(H' + $r + '="sg","Write the singular form on your keyboard - then press enter","Write the plural form on your keyboard - then press enter")'
    $ws.Range("I" + $r).Formula = $formula
}

# Move the live selection, matching where the author's cursor ended up.
$ws.Range("K7").Select()
